$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5465782294729991
$ws.Range("C2").Value = 0.130785156080016
$ws.Range("D2").Value = 0.03857308883073074
$ws.Range("F2").Value = 0.8229815943559302
$ws.Range("G2").Value = 0.002445606924888623
$ws.Range("I2").Value = 0.7804802698082511
$ws.Range("K2").Value = 0.3131678726726363
$ws.Range("L2").Value = 0.3022555150796791
$ws.Range("N2").Value = 1.581513752659225
$ws.Range("O2").Value = 2.889399287511651
$ws.Range("B3").Value = 0.5040314255320766
$ws.Range("C3").Value = 0.1303508687859534
$ws.Range("D3").Value = 0.03634294272045224
$ws.Range("F3").Value = 0.8224782102880326
$ws.Range("G3").Value = 0.002447861849060158
$ws.Range("I3").Value = 0.7853910922506415
$ws.Range("K3").Value = 0.2768638645424062
$ws.Range("L3").Value = 0.2909814548076071
$ws.Range("N3").Value = 1.59778547477941
$ws.Range("O3").Value = 2.900550683057574
$ws.Range("B4").Value = 0.4780530861884245
$ws.Range("C4").Value = 0.1300886977999873
$ws.Range("D4").Value = 0.03495945229420983
$ws.Range("F4").Value = 0.8226300606944434
$ws.Range("G4").Value = 0.002449321403288607
$ws.Range("I4").Value = 0.788780804884567
$ws.Range("K4").Value = 0.2545689930595643
$ws.Range("L4").Value = 0.284229878446439
$ws.Range("N4").Value = 1.608289682939263
$ws.Range("O4").Value = 2.908896192161137
$ws.Range("B5").Value = 0.4675039860581478
$ws.Range("C5").Value = 0.1299830086135749
$ws.Range("D5").Value = 0.03439213045037093
$ws.Range("F5").Value = 0.8228079177683654
$ws.Range("G5").Value = 0.002449935103208905
$ws.Range("I5").Value = 0.7902563337024269
$ws.Range("K5").Value = 0.2454831448119137
$ws.Range("L5").Value = 0.2815215766953116
$ws.Range("N5").Value = 1.612699314915353
$ws.Range("O5").Value = 2.912674005629569
$ws.Range("B6").Value = 0.4657545878745566
$ws.Range("C6").Value = 0.129965528915573
$ws.Range("D6").Value = 0.0342977141326628
$ws.Range("F6").Value = 0.8228444576376148
$ws.Range("G6").Value = 0.002450038151829465
$ws.Range("I6").Value = 0.7905070342363736
$ws.Range("K6").Value = 0.2439744311574827
$ws.Range("L6").Value = 0.2810744661091888
$ws.Range("N6").Value = 1.613439327421214
$ws.Range("O6").Value = 2.913324079687456
$ws.Range("B7").Value = 0.477910665472848
$ws.Range("C7").Value = 0.130087267765397
$ws.Range("D7").Value = 0.03495181548657911
$ws.Range("F7").Value = 0.8226319896561876
$ws.Range("G7").Value = 0.002449329603211813
$ws.Range("I7").Value = 0.7888003229499958
$ws.Range("K7").Value = 0.2544464595037397
$ws.Range("L7").Value = 0.2841931790394199
$ws.Range("N7").Value = 1.608348630120721
$ws.Range("O7").Value = 2.908945614632856
$ws.Range("B8").Value = 0.5318782722263222
$ws.Range("C8").Value = 0.1306344953206704
$ws.Range("D8").Value = 0.03780709039703822
$ws.Range("F8").Value = 0.8227124023369683
$ws.Range("G8").Value = 0.002446368885193167
$ws.Range("I8").Value = 0.782095825504161
$ws.Range("K8").Value = 0.3006513919282554
$ws.Range("L8").Value = 0.2983328290718816
$ws.Range("N8").Value = 1.587017711763343
$ws.Range("O8").Value = 2.892933352474017
$ws.Range("B9").Value = 0.6388392174477531
$ws.Range("C9").Value = 0.1317423517582981
$ws.Range("D9").Value = 0.04329294905335246
$ws.Range("F9").Value = 0.8265253456790234
$ws.Range("G9").Value = 0.002441155680199358
$ws.Range("I9").Value = 0.771918150273244
$ws.Range("K9").Value = 0.3912086156619523
$ws.Range("L9").Value = 0.3274133815234279
$ws.Range("N9").Value = 1.549259739661559
$ws.Range("O9").Value = 2.873419488581987
$ws.Range("B10").Value = 0.718087908014752
$ws.Range("C10").Value = 0.1325764580045856
$ws.Range("D10").Value = 0.04725348057240808
$ws.Range("F10").Value = 0.8315540409900208
$ws.Range("G10").Value = 0.002437683383364854
$ws.Range("I10").Value = 0.7662498691500446
$ws.Range("K10").Value = 0.4576920819700376
$ws.Range("L10").Value = 0.3496034065886278
$ws.Range("N10").Value = 1.523999074923352
$ws.Range("O10").Value = 2.866326930049979
$ws.Range("B11").Value = 0.7542796867443258
$ws.Range("C11").Value = 0.1329600843135736
$ws.Range("D11").Value = 0.04903988434727324
$ws.Range("F11").Value = 0.8343252597321751
$ws.Range("G11").Value = 0.002436180700956311
$ws.Range("I11").Value = 0.764063865556686
$ws.Range("K11").Value = 0.4879230335340594
$ws.Range("L11").Value = 0.3598773841203666
$ws.Range("N11").Value = 1.513045410113914
$ws.Range("O11").Value = 2.864673228757312
$ws.Range("B12").Value = 0.7680042812030763
$ws.Range("C12").Value = 0.1331059365256593
$ws.Range("D12").Value = 0.04971413240183864
$ws.Range("F12").Value = 0.8354441519115952
$ws.Range("G12").Value = 0.002435622673632851
$ws.Range("I12").Value = 0.7632925055429709
$ws.Range("K12").Value = 0.4993684380289096
$ws.Range("L12").Value = 0.363793650090301
$ws.Range("N12").Value = 1.508974839765262
$ws.Range("O12").Value = 2.86427310783705
$ws.Range("B13").Value = 0.765047586440204
$ws.Range("C13").Value = 0.133074499121264
$ws.Range("D13").Value = 0.04956902041800504
$ws.Range("F13").Value = 0.8352000886822921
$ws.Range("G13").Value = 0.00243574236602468
$ws.Range("I13").Value = 0.7634561220923288
$ws.Range("K13").Value = 0.496903581290411
$ws.Range("L13").Value = 0.3629490691744195
$ws.Range("N13").Value = 1.509848069210468
$ws.Range("O13").Value = 2.864349226121107
$ws.Range("B14").Value = 0.755408429777674
$ws.Range("C14").Value = 0.1329720721322545
$ws.Range("D14").Value = 0.049095399888337
$ws.Range("F14").Value = 0.8344159192651119
$ws.Range("G14").Value = 0.002436134571531079
$ws.Range("I14").Value = 0.7639992743933632
$ws.Range("K14").Value = 0.4888647052491706
$ws.Range("L14").Value = 0.360199062492768
$ws.Range("N14").Value = 1.512708971456423
$ws.Range("O14").Value = 2.864635779527021
$ws.Range("B15").Value = 0.7495066923171692
$ws.Range("C15").Value = 0.1329094077240853
$ws.Range("D15").Value = 0.04880500308928504
$ws.Range("F15").Value = 0.8339446410861342
$ws.Range("G15").Value = 0.002436376240252067
$ws.Range("I15").Value = 0.7643393196104427
$ws.Range("K15").Value = 0.4839403318907785
$ws.Range("L15").Value = 0.3585179524038438
$ws.Range("N15").Value = 1.514471430712623
$ws.Range("O15").Value = 2.864840744896355
$ws.Range("B16").Value = 0.71572546690237
$ws.Range("C16").Value = 0.1325514695946595
$ws.Range("D16").Value = 0.04713642524401251
$ws.Range("F16").Value = 0.8313826627162797
$ws.Range("G16").Value = 0.002437783129806389
$ws.Range("I16").Value = 0.7664006271712935
$ws.Range("K16").Value = 0.4557161134488297
$ws.Range("L16").Value = 0.3489355843060622
$ws.Range("N16").Value = 1.524725743323578
$ws.Range("O16").Value = 2.866466645908815
$ws.Range("B17").Value = 0.6950374024507937
$ws.Range("C17").Value = 0.1323329437069134
$ws.Range("D17").Value = 0.04610887678352782
$ws.Range("F17").Value = 0.8299348048998851
$ws.Range("G17").Value = 0.002438665867032563
$ws.Range("I17").Value = 0.7677656989190993
$ws.Range("K17").Value = 0.4383978065851295
$ws.Range("L17").Value = 0.3431030499017282
$ws.Range("N17").Value = 1.531154141407398
$ws.Range("O17").Value = 2.867866858742417
$ws.Range("B18").Value = 0.6831515053752071
$ws.Range("C18").Value = 0.1322076487555677
$ws.Range("D18").Value = 0.04551642334447337
$ws.Range("F18").Value = 0.8291475612151018
$ws.Range("G18").Value = 0.002439180833617161
$ws.Range("I18").Value = 0.7685877992380554
$ws.Range("K18").Value = 0.4284356124200883
$ws.Range("L18").Value = 0.3397652432480385
$ws.Range("N18").Value = 1.534902178611851
$ws.Range("O18").Value = 2.868820261314511
$ws.Range("B19").Value = 0.6791294570278694
$ws.Range("C19").Value = 0.1321652945777885
$ws.Range("D19").Value = 0.04531558327250451
$ws.Range("F19").Value = 0.8288888347999688
$ws.Range("G19").Value = 0.002439356437336806
$ws.Range("I19").Value = 0.7688724943753797
$ws.Range("K19").Value = 0.4250624025069669
$ws.Range("L19").Value = 0.3386380269726033
$ws.Range("N19").Value = 1.536179887542208
$ws.Range("O19").Value = 2.869168494330722
$ws.Range("B20").Value = 0.6972383079683766
$ws.Range("C20").Value = 0.1323561653948246
$ws.Range("D20").Value = 0.0462184097765217
$ws.Range("F20").Value = 0.8300842204784544
$ws.Range("G20").Value = 0.002438571149354261
$ws.Range("I20").Value = 0.7676165609002226
$ws.Range("K20").Value = 0.4402414943380109
$ws.Range("L20").Value = 0.3437221834199562
$ws.Range("N20").Value = 1.530464590970183
$ws.Range("O20").Value = 2.867702482817833
$ws.Range("B21").Value = 0.7582391590098894
$ws.Range("C21").Value = 0.13300214181249
$ws.Range("D21").Value = 0.04923457432175127
$ws.Range("F21").Value = 0.8346443632879073
$ws.Range("G21").Value = 0.002436019072886127
$ws.Range("I21").Value = 0.7638382059463282
$ws.Range("K21").Value = 0.4912259877329745
$ws.Range("L21").Value = 0.3610061084839202
$ws.Range("N21").Value = 1.511866555612089
$ws.Range("O21").Value = 2.864545476028752
$ws.Range("B22").Value = 0.7982203672794981
$ws.Range("C22").Value = 0.1334277031427504
$ws.Range("D22").Value = 0.05119283355607251
$ws.Range("F22").Value = 0.8380297053105465
$ws.Range("G22").Value = 0.002434415271027407
$ws.Range("I22").Value = 0.7616977561465816
$ws.Range("K22").Value = 0.5245330416841227
$ws.Range("L22").Value = 0.372452102599695
$ws.Range("N22").Value = 1.500162512730676
$ws.Range("O22").Value = 2.86380004491582
$ws.Range("B23").Value = 0.7768714933042986
$ws.Range("C23").Value = 0.1332002707898852
$ws.Range("D23").Value = 0.05014887111499178
$ws.Range("F23").Value = 0.8361858431757554
$ws.Range("G23").Value = 0.002435265398947365
$ws.Range("I23").Value = 0.7628100621520417
$ws.Range("K23").Value = 0.5067579343876787
$ws.Range("L23").Value = 0.3663294717631373
$ws.Range("N23").Value = 1.506367911736713
$ws.Range("O23").Value = 2.864077331165021
$ws.Range("B24").Value = 0.6962432537747532
$ws.Range("C24").Value = 0.1323456658143627
$ws.Range("D24").Value = 0.0461688952102719
$ws.Range("F24").Value = 0.8300165290554062
$ws.Range("G24").Value = 0.002438613947869376
$ws.Range("I24").Value = 0.767683870002827
$ws.Range("K24").Value = 0.4394079807199489
$ws.Range("L24").Value = 0.3434422251972222
$ws.Range("N24").Value = 1.530776173920721
$ws.Range("O24").Value = 2.867776334906807
$ws.Range("B25").Value = 0.6097849651521017
$ws.Range("C25").Value = 0.1314390331582302
$ws.Range("D25").Value = 0.04182109805526579
$ws.Range("F25").Value = 0.8251026256397722
$ws.Range("G25").Value = 0.002442502896867545
$ws.Range("I25").Value = 0.7743536682210994
$ws.Range("K25").Value = 0.3667177677061773
$ws.Range("L25").Value = 0.3194015098239902
$ws.Range("N25").Value = 1.55903858645513
$ws.Range("O25").Value = 2.877426200844724

Write-Output "done"